$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell E7 currently holds "...ser: 102" -> bump serial to 104
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 104"

# Cell C7 currently holds "...ser: 104" -> bump serial to 105 (new blog post)
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 105"
